$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.919.08'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '1.981.14'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.12'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.628'
$ws.Range("E6").Value = '  +1.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.89'
$ws.Range("E7").Value = '  +3.17%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.383'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0801'
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.90'
$ws.Range("E12").Value = '  +8.75%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.845'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.10'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = '2.276.79'
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.46'
$ws.Range("E16").Value = '  +3.77%  '
$ws.Range("D17").Value = '1.985.45'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = '36.860.18'
$ws.Range("E18").Value = '  +0.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.29'
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").Value = '0.0₃0861'
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.18'
$ws.Range("E21").Value = '  +2.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.86'
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.51'
$ws.Range("E24").Value = '  +3.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("E26").Value = '  +3.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.28'
$ws.Range("E27").Value = '  +0.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '163.53'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.54'
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.35'
$ws.Range("E30").Value = '  +17.66%  '
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.86'
$ws.Range("E32").Value = '  +3.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0620'
$ws.Range("E33").Value = '  +0.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.52'
$ws.Range("E34").Value = '  +5.71%  '
$ws.Range("B35").Value = 'BinanceUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.27'
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.56'
$ws.Range("E39").Value = '  -7.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0997'
$ws.Range("E40").Value = '  +1.27%  '
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("E42").Value = '  +0.99%  '
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.43'
$ws.Range("E44").Value = '  +1.55%  '
$ws.Range("D45").Value = '1.374.94'
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.10'
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.04'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.29'
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.00'
$ws.Range("E49").Value = '  +13.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.81'
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.15'
$ws.Range("E51").Value = '  +5.39%  '
